$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.024039885159386
$ws.Range("D2").Value = 1.025619098519711
$ws.Range("E2").Value = 1.024574542977473
$ws.Range("F2").Value = 1.034344367757102
$ws.Range("I2").Value = 1.030861816150708
$ws.Range("J2").Value = 1.029216730062773
$ws.Range("K2").Value = 1.028444139854746
$ws.Range("L2").Value = 1.027402644295856
$ws.Range("M2").Value = 1.037144111139605
$ws.Range("N2").Value = 1.030678335872274
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.024986486276833
$ws.Range("D3").Value = 1.026432725047172
$ws.Range("E3").Value = 1.025377938361743
$ws.Range("F3").Value = 1.035500464225372
$ws.Range("I3").Value = 1.030995613492662
$ws.Range("J3").Value = 1.029802082348963
$ws.Range("K3").Value = 1.029065090571623
$ws.Range("L3").Value = 1.02801316990532
$ws.Range("M3").Value = 1.038108451761614
$ws.Range("N3").Value = 1.031264519425849
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.025599315039097
$ws.Range("D4").Value = 1.026959783217269
$ws.Range("E4").Value = 1.025898442474207
$ws.Range("F4").Value = 1.036248957963323
$ws.Range("I4").Value = 1.031080776443965
$ws.Range("J4").Value = 1.030180566611908
$ws.Range("K4").Value = 1.029466826746031
$ws.Range("L4").Value = 1.028408224288188
$ws.Range("M4").Value = 1.038732293235059
$ws.Range("N4").Value = 1.031643541179859
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.025857022854793
$ws.Range("D5").Value = 1.027181497982094
$ws.Range("E5").Value = 1.026117417725073
$ws.Range("F5").Value = 1.036563725705226
$ws.Range("I5").Value = 1.031116240294867
$ws.Range("J5").Value = 1.030339614132392
$ws.Range("K5").Value = 1.029635701324173
$ws.Range("L5").Value = 1.028574305137942
$ws.Range("M5").Value = 1.03899451974487
$ws.Range("N5").Value = 1.031802814566065
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.025900297485617
$ws.Range("D6").Value = 1.027218733015214
$ws.Range("E6").Value = 1.026154193708059
$ws.Range("F6").Value = 1.036616582513765
$ws.Range("I6").Value = 1.031122174960792
$ws.Range("J6").Value = 1.030366314944576
$ws.Range("K6").Value = 1.029664055186884
$ws.Range("L6").Value = 1.028602190817994
$ws.Range("M6").Value = 1.039038546586736
$ws.Range("N6").Value = 1.031829553296465
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.025602758251268
$ws.Range("D7").Value = 1.026962745234031
$ws.Range("E7").Value = 1.025901367823358
$ws.Range("F7").Value = 1.036253163508604
$ws.Range("I7").Value = 1.031081251644864
$ws.Range("J7").Value = 1.030182692077911
$ws.Range("K7").Value = 1.029469083317616
$ws.Range("L7").Value = 1.028410443469748
$ws.Range("M7").Value = 1.038735797261978
$ws.Range("N7").Value = 1.031645669664268
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.024359727890756
$ws.Range("D8").Value = 1.025893945118831
$ws.Range("E8").Value = 1.024845918193154
$ws.Range("F8").Value = 1.034734988554676
$ws.Range("I8").Value = 1.030907325888465
$ws.Range("J8").Value = 1.029414609676055
$ws.Range("K8").Value = 1.028654005073875
$ws.Range("L8").Value = 1.027608973052063
$ws.Range("M8").Value = 1.03747004579265
$ws.Range("N8").Value = 1.030876496497302
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.022171789231076
$ws.Range("D9").Value = 1.024015130614591
$ws.Range("E9").Value = 1.022991137150457
$ws.Range("F9").Value = 1.032063004569645
$ws.Range("I9").Value = 1.030590043512438
$ws.Range("J9").Value = 1.028059055891224
$ws.Range("K9").Value = 1.027217306720418
$ws.Range("L9").Value = 1.026196743789224
$ws.Range("M9").Value = 1.035238479047147
$ws.Range("N9").Value = 1.029519017670612
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.020714843218652
$ws.Range("D10").Value = 1.022765710381692
$ws.Range("E10").Value = 1.021758082840325
$ws.Range("F10").Value = 1.030283856132507
$ws.Range("I10").Value = 1.030371282341967
$ws.Range("J10").Value = 1.027153988862461
$ws.Range("K10").Value = 1.026259269613481
$ws.Range("L10").Value = 1.025255349949011
$ws.Range("M10").Value = 1.033750002908995
$ws.Range("N10").Value = 1.028612665342855
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.02008437565198
$ws.Range("D11").Value = 1.022225451674446
$ws.Range("E11").Value = 1.021224992693616
$ws.Range("F11").Value = 1.029513979813175
$ws.Range("I11").Value = 1.03027484402769
$ws.Range("J11").Value = 1.026761770058727
$ws.Range("K11").Value = 1.025844382346254
$ws.Range("L11").Value = 1.024847748000185
$ws.Range("M11").Value = 1.033105295850557
$ws.Range("N11").Value = 1.028219889543432
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.019850252520885
$ws.Range("D12").Value = 1.022024889028964
$ws.Range("E12").Value = 1.021027105130656
$ws.Range("F12").Value = 1.02922808938136
$ws.Range("I12").Value = 1.030238765492769
$ws.Range("J12").Value = 1.026616035300025
$ws.Range("K12").Value = 1.025690267851332
$ws.Range("L12").Value = 1.024696351436395
$ws.Range("M12").Value = 1.032865794806303
$ws.Range("N12").Value = 1.028073947824657
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.019900469976376
$ws.Range("D13").Value = 1.022067905244748
$ws.Range("E13").Value = 1.021069546972217
$ws.Range("F13").Value = 1.029289410417672
$ws.Range("I13").Value = 1.030246516087453
$ws.Range("J13").Value = 1.02664729802711
$ws.Range("K13").Value = 1.0257233262334
$ws.Range("L13").Value = 1.024728826276934
$ws.Range("M13").Value = 1.032917169867054
$ws.Range("N13").Value = 1.0281052549484
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.020065021721628
$ws.Range("D14").Value = 1.022208870782364
$ws.Range("E14").Value = 1.021208632672324
$ws.Range("F14").Value = 1.029490346473535
$ws.Range("I14").Value = 1.030271867006474
$ws.Range("J14").Value = 1.026749724537767
$ws.Range("K14").Value = 1.02583164333025
$ws.Range("L14").Value = 1.024835233407582
$ws.Range("M14").Value = 1.03308549919369
$ws.Range("N14").Value = 1.028207826916451
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.020166415528633
$ws.Range("D15").Value = 1.022295739372604
$ws.Range("E15").Value = 1.021294344671615
$ws.Range("F15").Value = 1.029614159862042
$ws.Range("I15").Value = 1.0302874524931
$ws.Range("J15").Value = 1.026812826652775
$ws.Range("K15").Value = 1.025898380177919
$ws.Range("L15").Value = 1.024900795031468
$ws.Range("M15").Value = 1.033189208718045
$ws.Range("N15").Value = 1.028271018643699
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.020756693742031
$ws.Range("D16").Value = 1.022801581431596
$ws.Range("E16").Value = 1.021793479867641
$ws.Range("F16").Value = 1.030334960945479
$ws.Range("I16").Value = 1.030377646568252
$ws.Range("J16").Value = 1.027180012472646
$ws.Range("K16").Value = 1.026286803307341
$ws.Range("L16").Value = 1.025282401799533
$ws.Range("M16").Value = 1.033792786049031
$ws.Range("N16").Value = 1.028638725909551
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.021127067060154
$ws.Range("D17").Value = 1.023119083919991
$ws.Range("E17").Value = 1.022106797301689
$ws.Range("F17").Value = 1.03078723592689
$ws.Range("I17").Value = 1.030433764481227
$ws.Range("J17").Value = 1.02741025355279
$ws.Range("K17").Value = 1.026530437784669
$ws.Range("L17").Value = 1.025521781630646
$ws.Range("M17").Value = 1.034171344055468
$ws.Range("N17").Value = 1.028869293958438
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.021343138075916
$ws.Range("D18").Value = 1.023304349899348
$ws.Range("E18").Value = 1.022289630119619
$ws.Range("F18").Value = 1.031051089174427
$ws.Range("I18").Value = 1.030466331779849
$ws.Range("J18").Value = 1.027544518474505
$ws.Range("K18").Value = 1.026672540693903
$ws.Range("L18").Value = 1.02566141046633
$ws.Range("M18").Value = 1.034392132488553
$ws.Range("N18").Value = 1.029003749551743
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.021416819226461
$ws.Range("D19").Value = 1.023367533016768
$ws.Range("E19").Value = 1.022351984932342
$ws.Range("F19").Value = 1.031141064614558
$ws.Range("I19").Value = 1.030477408333477
$ws.Range("J19").Value = 1.027590294097778
$ws.Range("K19").Value = 1.026720993249568
$ws.Range("L19").Value = 1.025709020713578
$ws.Range("M19").Value = 1.034467412578434
$ws.Range("N19").Value = 1.029049590181651
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.021087325519392
$ws.Range("D20").Value = 1.023085011437475
$ws.Range("E20").Value = 1.022073173017454
$ws.Range("F20").Value = 1.030738706044831
$ws.Range("I20").Value = 1.030427760658627
$ws.Range("J20").Value = 1.02738555404532
$ws.Range("K20").Value = 1.026504298618483
$ws.Range("L20").Value = 1.025496098162567
$ws.Range("M20").Value = 1.034130730230855
$ws.Range("N20").Value = 1.028844559374834
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.020016563656862
$ws.Range("D21").Value = 1.02216735680369
$ws.Range("E21").Value = 1.021167671917434
$ws.Range("F21").Value = 1.029431173719434
$ws.Range("I21").Value = 1.030264408885728
$ws.Range("J21").Value = 1.026719563781765
$ws.Range("K21").Value = 1.025799746829142
$ws.Range("L21").Value = 1.024803899016931
$ws.Range("M21").Value = 1.03303593117946
$ws.Range("N21").Value = 1.028177623328715
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.019343684031876
$ws.Range("D22").Value = 1.021591048065582
$ws.Range("E22").Value = 1.02059907611849
$ws.Range("F22").Value = 1.028609515351019
$ws.Range("I22").Value = 1.030160216098297
$ws.Range("J22").Value = 1.026300556332925
$ws.Range("K22").Value = 1.025356727658614
$ws.Range("L22").Value = 1.024368714953029
$ws.Range("M22").Value = 1.032347424756798
$ws.Range("N22").Value = 1.027758020841237
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.01970035681125
$ws.Range("D23").Value = 1.021896497594479
$ws.Range("E23").Value = 1.020900430079156
$ws.Range("F23").Value = 1.029045050447016
$ws.Range("I23").Value = 1.03021559150378
$ws.Range("J23").Value = 1.026522705758445
$ws.Range("K23").Value = 1.025591583995754
$ws.Range("L23").Value = 1.024599411251517
$ws.Range("M23").Value = 1.032712430480637
$ws.Range("N23").Value = 1.027980485744421
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.021105282891948
$ws.Range("D24").Value = 1.023100407102839
$ws.Range("E24").Value = 1.022088366136875
$ws.Range("F24").Value = 1.030760634458288
$ws.Range("I24").Value = 1.030430474038488
$ws.Range("J24").Value = 1.027396714784776
$ws.Range("K24").Value = 1.026516109797899
$ws.Range("L24").Value = 1.025507703408179
$ws.Range("M24").Value = 1.034149081925602
$ws.Range("N24").Value = 1.02885573596382
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.022737130649063
$ws.Range("D25").Value = 1.024500303546244
$ws.Range("E25").Value = 1.023470036432208
$ws.Range("F25").Value = 1.032753393243794
$ws.Range("I25").Value = 1.030673346376403
$ws.Range("J25").Value = 1.028409742737421
$ws.Range("K25").Value = 1.027588773211445
$ws.Range("L25").Value = 1.026561826235443
$ws.Range("M25").Value = 1.035815528029309
$ws.Range("N25").Value = 1.029870202532355
